$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")

# Add the new worksheet "Hoja1" positioned after Sheet1
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "Hoja1"
$ws.Columns.Item(1).ColumnWidth = 12.125

# --- Phase 1: seed each unique text value once, in original authoring order,
# so the shared-strings table is built up in the same sequence as the source file ---
$ws.Range("A1").Value = "Mat original"
$ws.Range("A5").Value = "mask inpaint"
$ws.Range("A12").Value = "mat extended"
$ws.Range("T20").Value = "b"
$ws.Range("B12").Value = "b1"
$ws.Range("C12").Value = "b2"
$ws.Range("D12").Value = "b3"
$ws.Range("E12").Value = "b4"
$ws.Range("B13").Value = "b5"
$ws.Range("B14").Value = "b7"
$ws.Range("B15").Value = "b9"
$ws.Range("C15").Value = "b10"
$ws.Range("D15").Value = "b11"
$ws.Range("E15").Value = "b12"
$ws.Range("E13").Value = "b6"
$ws.Range("E14").Value = "b8"
$ws.Range("R20").Value = "x"
$ws.Range("A19").Value = "mat a"

# --- Phase 2: fill in the remaining cells (numbers + repeated text) ---
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 4
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 2
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 4
$ws.Range("B19").Value = "b1"
$ws.Range("C19").Value = "b5"
$ws.Range("D19").Value = "b7"
$ws.Range("E19").Value = "b9"
$ws.Range("F19").Value = "b2"
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 3
$ws.Range("I19").Value = "b10"
$ws.Range("J19").Value = "b3"
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 3
$ws.Range("M19").Value = "b11"
$ws.Range("N19").Value = "b4"
$ws.Range("O19").Value = "b6"
$ws.Range("P19").Value = "b8"
$ws.Range("Q19").Value = "b12"
$ws.Range("A20").Value = "b1"
$ws.Range("B20").Value = 2
$ws.Range("B20").Font.Underline = $true
$ws.Range("C20").Value = -1
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = -1
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0
$ws.Range("U20").Value = 0
$ws.Range("A21").Value = "b5"
$ws.Range("U21").Value = "b5"
$ws.Range("A22").Value = "b7"
$ws.Range("U22").Value = "b7"
$ws.Range("A23").Value = "b9"
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = -1
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0
$ws.Range("U23").Value = 0
$ws.Range("A24").Value = "b2"
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = -1
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0
$ws.Range("U24").Value = 0
$ws.Range("A25").Value = 1
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = -1
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = -1
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = -1
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = -1
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0
$ws.Range("U25").Value = 0
$ws.Range("A26").Value = 3
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 0
$ws.Range("U26").Value = 3
$ws.Range("A27").Value = "b10"
$ws.Range("U27").Value = "b10"
$ws.Range("A28").Value = "b3"
$ws.Range("U28").Value = "b3"
$ws.Range("A29").Value = 2
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0
$ws.Range("U29").Value = 2
$ws.Range("A30").Value = 4
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 1
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("Q30").Value = 0
$ws.Range("U30").Value = 4
$ws.Range("A31").Value = "b11"
$ws.Range("U31").Value = "b11"
$ws.Range("A32").Value = "b4"
$ws.Range("U32").Value = "b4"
$ws.Range("A33").Value = "b6"
$ws.Range("U33").Value = "b6"
$ws.Range("A34").Value = "b8"
$ws.Range("U34").Value = "b8"
$ws.Range("A35").Value = "b12"
$ws.Range("U35").Value = "b12"

# --- View/selection state ---
$ws.Range("G25").Select() | Out-Null

$sheet1.Select() | Out-Null
$sheet1.Range("BZ1:BZ8").Select() | Out-Null

# --- Page setup (paper size / orientation) on both sheets ---
$sheet1.PageSetup.PaperSize = 9
$sheet1.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Hoja1 is the tab that was active/selected when the file was saved
$ws.Select() | Out-Null
